$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$trackNums = @{
    2  = "320018722164"
    3  = "320018722370"
    4  = "320018722429"
    5  = "320018722440"
    6  = "320018722484"
    7  = "320018722510"
    8  = "320018722543"
    9  = "320018722565"
    10 = "320018722602"
    11 = "320018722646"
    12 = "320018722680"
    13 = "320018722716"
    14 = "320018722749"
    15 = "320018722782"
    16 = "320018722819"
    17 = "320018722830"
    18 = "320018722874"
    19 = "320018722896"
    20 = "320018722922"
    21 = "320018722944"
    22 = "320018722977"
}

$mirrorRows = @(5, 6, 7, 13, 14, 15, 16, 17)

foreach ($row in $trackNums.Keys) {
    $value = $trackNums[$row]
    $cellC = $ws.Range("C$row")
    $cellC.NumberFormat = "@"
    $cellC.Value = $value
    $cellC.Style = "Normal"
    if ($mirrorRows -contains $row) {
        $cellD = $ws.Range("D$row")
        $cellD.NumberFormat = "@"
        $cellD.Value = $value
        $cellD.Style = "Normal"
    }
}
